$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.135675479443519
$ws.Cells.Item(2, 3).Value = 0.473842788633127
$ws.Cells.Item(3, 2).Value = 0.135198094040841
$ws.Cells.Item(3, 3).Value = 0.460791351757007
$ws.Cells.Item(4, 2).Value = 0.0581384868268255
$ws.Cells.Item(4, 3).Value = 0.279988761045496
$ws.Cells.Item(5, 2).Value = 0.215338968844249
$ws.Cells.Item(5, 3).Value = 0.632623940737174
$ws.Cells.Item(6, 2).Value = 0.212810850633746
$ws.Cells.Item(6, 3).Value = 0.631192968656138
$ws.Cells.Item(7, 2).Value = 0.124448274891151
$ws.Cells.Item(7, 3).Value = 0.490201736805172
$ws.Cells.Item(8, 2).Value = 0.191256596210678
$ws.Cells.Item(8, 3).Value = 0.774046143519458
$ws.Cells.Item(9, 2).Value = 0.210623869643663
$ws.Cells.Item(9, 3).Value = 0.792234563204733
$ws.Cells.Item(10, 2).Value = 0.229838728344693
$ws.Cells.Item(10, 3).Value = 0.659700510490542
$ws.Cells.Item(11, 2).Value = 0.187851850948922
$ws.Cells.Item(11, 3).Value = 0.729332548739581
$ws.Cells.Item(12, 2).Value = 0.158985275362944
$ws.Cells.Item(12, 3).Value = 0.481205105399773
$ws.Cells.Item(13, 2).Value = 0.135650403714869
$ws.Cells.Item(13, 3).Value = 0.929001763996532
$ws.Cells.Item(14, 2).Value = 0.199890021605445
$ws.Cells.Item(14, 3).Value = 0.633568644801997
$ws.Cells.Item(15, 2).Value = 0.182031513646177
$ws.Cells.Item(15, 3).Value = 0.829540995478474
$ws.Cells.Item(16, 2).Value = 0.222696039185922
$ws.Cells.Item(16, 3).Value = 0.675641435518982
$ws.Cells.Item(17, 2).Value = 0.198484683707293
$ws.Cells.Item(17, 3).Value = 0.585302236064272
$ws.Cells.Item(18, 2).Value = 0.249156488427116
$ws.Cells.Item(18, 3).Value = 0.823820104523769
$ws.Cells.Item(19, 2).Value = 0.200572948437639
$ws.Cells.Item(19, 3).Value = 0.659159066760244
$ws.Cells.Item(20, 2).Value = 0.107189625075493
$ws.Cells.Item(20, 3).Value = 0.439344367434619
$ws.Cells.Item(21, 2).Value = 0.0940485505915715
$ws.Cells.Item(21, 3).Value = 0.463658576529257
$ws.Cells.Item(22, 2).Value = 0.162555350705141
$ws.Cells.Item(22, 3).Value = 0.515055030560506
$ws.Cells.Item(23, 2).Value = 0.233568655736157
$ws.Cells.Item(23, 3).Value = 0.749403991617224
$ws.Cells.Item(24, 2).Value = 0.145723933270154
$ws.Cells.Item(24, 3).Value = 0.576169097415042
$ws.Cells.Item(25, 2).Value = 0.108861627629659
$ws.Cells.Item(25, 3).Value = 0.385076627920173
$ws.Cells.Item(26, 2).Value = 0.266269000004038
$ws.Cells.Item(26, 3).Value = 0.863498245552782
$ws.Cells.Item(27, 2).Value = 0.273581980872965
$ws.Cells.Item(27, 3).Value = 0.641524963502602
$ws.Cells.Item(28, 2).Value = 0.227149154917128
$ws.Cells.Item(28, 3).Value = 0.861163482610652
$ws.Cells.Item(29, 2).Value = 0.12165853248952
$ws.Cells.Item(29, 3).Value = 0.844209360580092
$ws.Cells.Item(30, 2).Value = 0.118353035793507
$ws.Cells.Item(30, 3).Value = 0.953956548204336
$ws.Cells.Item(31, 2).Value = 0.137988604712299
$ws.Cells.Item(31, 3).Value = 0.961908741558781
$ws.Cells.Item(32, 2).Value = 0.185996675990196
$ws.Cells.Item(32, 3).Value = 0.850898837610772
$ws.Cells.Item(33, 2).Value = 0.109601031603587
$ws.Cells.Item(33, 3).Value = 0.979626847934591
$ws.Cells.Item(34, 2).Value = 0.184657162511518
$ws.Cells.Item(34, 3).Value = 0.957876693509183
